$wb = $excel.ActiveWorkbook

# Rename Sheet1 -> 接口, Sheet2 -> 登录
$wsApi = $wb.Worksheets.Item("Sheet1")
$wsApi.Name = "接口"

$wsLogin = $wb.Worksheets.Item("Sheet2")
$wsLogin.Name = "登录"

# Add the new "用户登录" (user login) record to the 接口 sheet
$wsApi.Cells.Item(2, 1).Value = 1
$wsApi.Cells.Item(2, 2).Value = "用户登录"
$wsApi.Cells.Item(2, 3).Value = "POST"
$wsApi.Cells.Item(2, 4).Value = "/account/login"
$wsApi.Cells.Item(2, 5).Value = "登录"

# Widen column D to fit the new URL text
$wsApi.Columns.Item(4).AutoFit()

# Update selection on the 接口 sheet
$wsApi.Range("F8").Select()

# Make the 登录 sheet the active tab
$wsLogin.Activate()
